# Auto-generated edit script: apply numeric corrections to the
# "Excalibur_Profits" per-craft profitability sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Values come from a scheduled data refresh; only
# specific H..N cells on specific rows change (no structural changes).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 3555.2222
$ws.Range("I76").Value = 3555.2222
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3555.2222
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3240.2222
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 3555.2222
$ws.Range("I79").Value = 3555.2222
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3555.2222
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2463.2222
$ws.Range("N79").ClearContents()
# Row 92
$ws.Range("H92").Value = 441.375
$ws.Range("I92").Value = 441.375
$ws.Range("K92").Value = 441.375
$ws.Range("M92").Value = 806.625
# Row 95
$ws.Range("H95").Value = 80467.75
$ws.Range("J95").Value = 80467.75
$ws.Range("L95").Value = 80467.75
$ws.Range("N95").Value = -85959.75
# Row 138
$ws.Range("H138").Value = 2069.121
$ws.Range("I138").Value = 1386.579
$ws.Range("J138").Value = 2995.4285
$ws.Range("K138").Value = 4159.737
$ws.Range("L138").Value = 8986.2855
$ws.Range("M138").Value = 980.2629999999999
$ws.Range("N138").Value = -19266.2855
# Row 141
$ws.Range("H141").Value = 3250
$ws.Range("I141").Value = 3000
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 9000
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -3820
$ws.Range("N141").Value = -22360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5379131.5
$ws.Range("I32").Value = 5558417
$ws.Range("K32").Value = 5558417
$ws.Range("M32").Value = -5558130
# Row 97
$ws.Range("H97").Value = 1215.909
$ws.Range("I97").Value = 1123.0714
$ws.Range("J97").Value = 1378.375
$ws.Range("K97").Value = 1123.0714
$ws.Range("L97").Value = 1378.375
$ws.Range("M97").Value = -627.0714
$ws.Range("N97").Value = -2370.375
# Row 102
$ws.Range("H102").Value = 36123
$ws.Range("I102").Value = 36123
$ws.Range("K102").Value = 36123
$ws.Range("M102").Value = -34501
# Row 122
$ws.Range("H122").Value = 2353.9656
$ws.Range("I122").Value = 2277.3215
$ws.Range("K122").Value = 6831.9645
$ws.Range("M122").Value = -4381.9645
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1452
$ws.Range("I94").Value = 1345.8572
$ws.Range("J94").Value = 1699.6666
$ws.Range("K94").Value = 1345.8572
$ws.Range("L94").Value = 1699.6666
$ws.Range("M94").Value = -894.8571999999999
$ws.Range("N94").Value = -2601.6666
# Row 99
$ws.Range("H99").Value = 3414
$ws.Range("I99").Value = 3816.3333
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 3816.3333
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -2318.3333
$ws.Range("N99").Value = -3996
# Row 105
$ws.Range("H105").Value = 1727.6
$ws.Range("I105").Value = 1634.5
$ws.Range("K105").Value = 1634.5
$ws.Range("M105").Value = 112.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 175976.3
$ws.Range("I31").Value = 285227.72
$ws.Range("J31").Value = 42446.777
$ws.Range("K31").Value = 285227.72
$ws.Range("L31").Value = 42446.777
$ws.Range("M31").Value = -284932.72
$ws.Range("N31").Value = -43036.777
# Row 34
$ws.Range("H34").Value = 175976.3
$ws.Range("I34").Value = 285227.72
$ws.Range("J34").Value = 42446.777
$ws.Range("K34").Value = 285227.72
$ws.Range("L34").Value = 42446.777
$ws.Range("M34").Value = -285025.72
$ws.Range("N34").Value = -42850.777
# Row 86
$ws.Range("H86").Value = 99982.48
$ws.Range("I86").Value = 4614.4546
$ws.Range("K86").Value = 4614.4546
$ws.Range("M86").Value = -3491.4546
# Row 89
$ws.Range("H89").Value = 99982.48
$ws.Range("I89").Value = 4614.4546
$ws.Range("K89").Value = 23072.273
$ws.Range("M89").Value = -17456.273
# Row 99
$ws.Range("H99").Value = 9887.166999999999
$ws.Range("I99").Value = 1569.5
$ws.Range("J99").Value = 14046
$ws.Range("K99").Value = 1569.5
$ws.Range("L99").Value = 14046
$ws.Range("M99").Value = -71.5
$ws.Range("N99").Value = -17042
# Row 126
$ws.Range("H126").Value = 9887.166999999999
$ws.Range("I126").Value = 1569.5
$ws.Range("J126").Value = 14046
$ws.Range("K126").Value = 4708.5
$ws.Range("L126").Value = 42138
$ws.Range("M126").Value = -2238.5
$ws.Range("N126").Value = -47078
# Row 134
$ws.Range("H134").Value = 7221.773
$ws.Range("I134").Value = 7470.4287
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 22411.2861
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -19876.2861
$ws.Range("N134").Value = -11070

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 133
$ws.Range("H133").Value = 7139
$ws.Range("I133").Value = 5665
$ws.Range("J133").Value = 8613
$ws.Range("K133").Value = 16995
$ws.Range("L133").Value = 25839
$ws.Range("M133").Value = -11935
$ws.Range("N133").Value = -35959

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 69
$ws.Range("H69").Value = 31577.834
$ws.Range("J69").Value = 31577.834
$ws.Range("L69").Value = 31577.834
$ws.Range("N69").Value = -33075.834
# Row 70
$ws.Range("H70").Value = 11364.9375
$ws.Range("J70").Value = 5170.75
$ws.Range("L70").Value = 5170.75
$ws.Range("N70").Value = -5710.75
# Row 72
$ws.Range("H72").Value = 31577.834
$ws.Range("J72").Value = 31577.834
$ws.Range("L72").Value = 94733.50199999999
$ws.Range("N72").Value = -102221.502
# Row 73
$ws.Range("H73").Value = 11364.9375
$ws.Range("J73").Value = 5170.75
$ws.Range("L73").Value = 5170.75
$ws.Range("N73").Value = -7042.75
# Row 97
$ws.Range("H97").Value = 1557.7587
$ws.Range("I97").Value = 1315.44
$ws.Range("J97").Value = 3072.25
$ws.Range("K97").Value = 1315.44
$ws.Range("L97").Value = 3072.25
$ws.Range("M97").Value = -819.4400000000001
$ws.Range("N97").Value = -4064.25
# Row 102
$ws.Range("H102").Value = 5565.8
$ws.Range("I102").Value = 5331.9062
$ws.Range("K102").Value = 5331.9062
$ws.Range("M102").Value = -3709.9062
# Row 132
$ws.Range("H132").Value = 21092338
$ws.Range("I132").Value = 31632288
$ws.Range("J132").Value = 12440.8125
$ws.Range("K132").Value = 94896864
$ws.Range("L132").Value = 37322.4375
$ws.Range("M132").Value = -94894334
$ws.Range("N132").Value = -42382.4375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1720.6666
$ws.Range("I68").Value = 1325
$ws.Range("J68").Value = 2314.1667
$ws.Range("K68").Value = 1325
$ws.Range("L68").Value = 2314.1667
$ws.Range("M68").Value = -576
$ws.Range("N68").Value = -3812.1667
# Row 71
$ws.Range("H71").Value = 1720.6666
$ws.Range("I71").Value = 1325
$ws.Range("J71").Value = 2314.1667
$ws.Range("K71").Value = 6625
$ws.Range("L71").Value = 11570.8335
$ws.Range("M71").Value = -2881
$ws.Range("N71").Value = -19058.8335
# Row 74
$ws.Range("H74").Value = 45806.43
$ws.Range("I74").Value = 19998.5
$ws.Range("K74").Value = 19998.5
$ws.Range("M74").Value = -19000.5
# Row 77
$ws.Range("H77").Value = 45806.43
$ws.Range("I77").Value = 19998.5
$ws.Range("K77").Value = 59995.5
$ws.Range("M77").Value = -55003.5
# Row 82
$ws.Range("H82").Value = 1866.6666
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 2300
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 2300
$ws.Range("M82").Value = -639
$ws.Range("N82").Value = -3022
# Row 85
$ws.Range("H85").Value = 1866.6666
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 2300
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 2300
$ws.Range("M85").Value = 248
$ws.Range("N85").Value = -4796
# Row 96
$ws.Range("H96").Value = 99750
$ws.Range("J96").Value = 89666.664
$ws.Range("L96").Value = 89666.664
$ws.Range("N96").Value = -95158.664
# Row 100
$ws.Range("H100").Value = 27251
$ws.Range("J100").Value = 100004
$ws.Range("L100").Value = 100004
$ws.Range("N100").Value = -101086
# Row 132
$ws.Range("H132").Value = 1744918.8
$ws.Range("I132").Value = 2679875.5
$ws.Range("K132").Value = 8039626.5
$ws.Range("M132").Value = -8037096.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 37514
$ws.Range("I136").Value = 32167.334
$ws.Range("J136").Value = 48207.332
$ws.Range("K136").Value = 96502.00199999999
$ws.Range("L136").Value = 144621.996
$ws.Range("M136").Value = -149721.996
